# Open new account test case: update the test "Username" value on the
# "Sign up" sheet and make that sheet/cell the active selection instead
# of the "URL" sheet.

$wb = $excel.ActiveWorkbook
$signUp = $wb.Worksheets.Item("Sign up")

# Update the username test data used for the "open new account" test.
$signUp.Range("I2").Value = "madboly5"

# Make "Sign up" the active sheet/tab and select the edited cell.
$signUp.Activate()
$signUp.Range("I2").Select()
